$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text: B1 shared string "value" -> "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Propagate A2's date style/number-format down through A3:A22 before writing values,
# so every date cell keeps the same border/font/alignment/number-format as A2.
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Date (A) / value (B) pairs for rows 2..22
$dates = @(38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657, 46022)
$values = @($null, 0.8557439673732903, 1.788430953138542, 1.807765267947059, 1.267704211901699, 0.815159612280536, 2.321967807433256, 1.090188641041823, 0.0688236519329477, 1.078691045907165, 2.724993727165903, 3.366078187926935, 1.625351934832997, 0.8215983724355613, 2.013614902241634, 2.422659863072885, 2.094164808525223, 0.6457774251209525, 1.355088394421644, 0.9622045007620983, $null)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    if ($null -ne $values[$i]) {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    } else {
        $ws.Cells.Item($row, 2).Value = $null
    }
}
